# "update shop and ball"
#  - add a new "ball" sheet (unlockable ball/colour shop data) right after "player"
#  - add a "locked" column to the "player" sheet

$wb = $excel.ActiveWorkbook

$playerSheet = $wb.Worksheets.Item("player")

# ---------------------------------------------------------------------------
# 1. player sheet: add column G "locked" / "解锁状态"
# ---------------------------------------------------------------------------
$playerSheet.Range("G1").Value = "locked"
$playerSheet.Range("G2").Value = "解锁状态"
$playerSheet.Range("G3").Value = "c"

for ($row = 4; $row -le 34; $row++) {
    $playerSheet.Cells.Item($row, 7).Value = 0
}

[void]$playerSheet.Range("A1:E3").Select()

# ---------------------------------------------------------------------------
# 2. Insert new "ball" sheet right after "player"
# ---------------------------------------------------------------------------
$ballSheet = $wb.Worksheets.Add($null, $playerSheet)
$ballSheet.Name = "ball"

# headers
$ballSheet.Range("A1").Value = "id"
$ballSheet.Range("B1").Value = "color"
$ballSheet.Range("C1").Value = "locked"
$ballSheet.Range("D1").Value = "index"

$ballSheet.Range("A2").Value = "编号"
$ballSheet.Range("B2").Value = "颜色"
$ballSheet.Range("C2").Value = "状态"
$ballSheet.Range("D2").Value = "序列"

$ballSheet.Range("A3").Value = "c"
$ballSheet.Range("B3").Value = "c"
$ballSheet.Range("C3").Value = "c"
$ballSheet.Range("D3").Value = "c"

$colors = @("#ffffff", "#ffea00", "#fe0000", "#6500ff", "#ff5a00", "#ff009a", "#06ff00", "#0075ff")

for ($i = 0; $i -lt $colors.Length; $i++) {
    $r = 4 + $i
    $ballSheet.Cells.Item($r, 1).Value = $i + 1
    $ballSheet.Cells.Item($r, 2).Value = $colors[$i]
    $ballSheet.Cells.Item($r, 3).Value = 0
    $ballSheet.Cells.Item($r, 4).Value = $i + 1
}

[void]$ballSheet.Range("E11").Select()
[void]$ballSheet.Activate()
